$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Manually shift rows 12-15 down to 13-16 to make room for a new "Jurisdiction" row at 12.
# (Using Range copies instead of Rows.Insert() avoids an extra/unused style entry being
# minted in styles.xml for the freshly inserted blank row.)
for ($r = 15; $r -ge 12; $r--) {
    $dst = $r + 1
    $ws1.Range("A$dst").Value = $ws1.Range("A$r").Value2
    $ws1.Range("B$dst").Value = $ws1.Range("B$r").Value2
}

# Row 16 is brand new (beyond the old A1:B15 range) so it has no style yet - copy the
# formatting from row 15 (which already carries the correct style) onto it.
$ws1.Range("A15:B15").Copy()
$ws1.Range("A16:B16").PasteSpecial(-4122)
$ws1.Range("A16").Value = "Immutable"
$ws1.Range("B16").Value = "BooleanType[null]"

# New row 12: Jurisdiction / (blank)
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""

# Version 0.1.6 -> 0.1.7
$ws1.Range("B3").Value = "0.1.7"
# Status active -> draft
$ws1.Range("B6").Value = "draft"
# Date updated
$ws1.Range("B8").Value = "2024-08-27T12:23:18-05:00"
# Contact (row 10) display text updated
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
# Contact (row 11) - second contact, was a duplicate, now a real person
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"
